# Add a "Date Time" column (I) to TestSheet3:
#  - header "Date Time" in I1 (same header style as the other headers)
#  - a sample datetime value in I2, formatted as yyyy/mm/dd hh:mm
#  - widen column I to fit the new content
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSheet3")

# Header cell: value + copy the formatting (bold font / bottom border) from H1
$ws.Range("I1").Value = "Date Time"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data cell: a date-time serial value with a custom date/time number format
$ws.Range("I2").Value = 43102.916666666664
$ws.Range("I2").NumberFormat = "yyyy/mm/dd\ hh:mm"

# Resize the new column so the date/time values are fully visible
$ws.Columns.Item(9).ColumnWidth = 19.5703125

# Keep the sheet's selection in sync with the new used range (A1:I3)
[void]$ws.Range("A1:I3").Select()

Write-Output "Added Date Time column to TestSheet3"
